# Avances del proyecto final de econometría
#
# The workbook has three sheets: "data", "var", "diccionario".
# Rename the "var" sheet to "vars" (the only substantive content change
# in this revision).
$wb = $excel.ActiveWorkbook

$ws = $null
try {
    $ws = $wb.Worksheets.Item("var")
} catch {
    $ws = $null
}

if ($ws -eq $null) {
    # Fallback: find it by name among the worksheets.
    foreach ($sheet in $wb.Worksheets) {
        if ($sheet.Name -eq "var") {
            $ws = $sheet
        }
    }
}

if ($ws -ne $null) {
    $ws.Name = "vars"
}
